$wb = $excel.ActiveWorkbook

# ---------- Sheet1 : add new reference rows 40-56 ----------
$ws1 = $wb.Worksheets.Item("Sheet1")

# shared text blocks used by the new Sheet1 rows
$s156 = @'
Swing
'@

$s157 = @'
Basic example
'@

$s186 = @'
package com.bcm.app.ui;
import javax.swing.JFrame;
import javax.swing.JLabel;
import javax.swing.JDialog;
import javax.swing.JButton;
import javax.swing.JCheckBox;
import javax.swing.SwingConstants;
import java.awt.event.ActionListener;
import java.awt.event.ActionEvent;
import java.awt.event.ItemListener;
import java.awt.event.ItemEvent;
import java.awt.EventQueue;
import java.awt.Font;
import java.sql.*;
import org.sqlite.SQLiteConfig;
import org.sqlite.SQLiteDataSource;
import com.bcm.app.engine.SendSMSJob;
public class SendSMSUI implements ActionListener, ItemListener {
    private SendSMSJob mJob;
    /* Main Frame properties */
    private JFrame mMainFrame;
    private JLabel mLastSentTimeTagLabel;
    private JLabel mLastSentTimeLabel;
    private JLabel mJobStatusTagLabel;
    private JLabel mJobStatusLabel;
    private JLabel mMomentLabel;
    private JButton mStartButton;
    private JButton mStopButton;
    private JButton mExportLogButton;
    private JButton mClearLogButton;
    private JButton mLoadConfigButton;
    private JButton mExitButton;
    /* Load Config Frame properties */
    private JFrame mLoadConfigFrame;
    private JLabel mFtpAddressTagLabel;
    private JLabel mFtpAddressLabel;
    // private JLabel mFtpPortTagLabel;
    // private JLabel mFtpPortLabel;
    // private JLabel mFtpUserTagLabel;
    // private JLabel mFtpUserLabel;
    // private JLabel mFtpPasswrodTagLabel;
    // private JLabel mFtpPasswordLabel;
    // more label to be added...
    private JCheckBox mFetchConfigCheckBox;
    private JButton mLoadConfirmButton;
    private JButton mLoadCancelButton;
    /**
     * Launch the application.
     */
    public static void main(String[] args) {
        EventQueue.invokeLater(new Runnable() {
            public void run() {
                try {
                    SendSMSUI window = new SendSMSUI();
                    window.mMainFrame.setVisible(true);
                } catch (Exception e) {
                    e.printStackTrace();
                }
            }
        });
    }
    /**
     * Create the application.
     */
    public SendSMSUI() {          
        mJob = new SendSMSJob();
        mJob.setProperties("config.properties");
        initialize();
    }
    /**
     * Initialize all Frames.
     */
    private void initialize (){
        initializeMain();
        initializeLoadConfig();
    }
    /**
     * Init mMainFrame
     */
    private void initializeMain() {
        mMainFrame = new JFrame();
        mMainFrame.setBounds(100, 100, 360, 360);
        mMainFrame.setDefaultCloseOperation(JFrame.EXIT_ON_CLOSE);
        mMainFrame.getContentPane().setLayout(null);
        /*Last sent time */
        mLastSentTimeTagLabel = new JLabel("Last Sent Time: ");
        mLastSentTimeTagLabel.setBounds(26, 32, 100, 16);
        mMainFrame.getContentPane().add(mLastSentTimeTagLabel);
        mLastSentTimeLabel = new JLabel("[...]");
        mLastSentTimeLabel.setBounds(136, 33, 175, 14);
        mMainFrame.getContentPane().add(mLastSentTimeLabel);
        /* Job status */
        mJobStatusTagLabel = new JLabel("Job Status: ");
        mJobStatusTagLabel.setBounds(26, 59, 100, 14);
        mMainFrame.getContentPane().add(mJobStatusTagLabel);
        mJobStatusLabel = new JLabel("[...]");
        mJobStatusLabel.setBounds(136, 59, 175, 14);
        mMainFrame.getContentPane().add(mJobStatusLabel);
        /* Moment status */
        mMomentLabel = new JLabel("MOMENT");
        mMomentLabel.setHorizontalAlignment(SwingConstants.CENTER);
        mMomentLabel.setFont(new Font("Tahoma", Font.PLAIN, 25));
        mMomentLabel.setBounds(0, 84, 344, 48);
        mMainFrame.getContentPane().add(mMomentLabel);
        /* Button: Start and stop */
        mStartButton = new JButton("Start Send SMS");
        mStartButton.setBounds(26, 154, 125, 25);
        mStartButton.addActionListener(this);
        mMainFrame.getContentPane().add(mStartButton);
        mStopButton = new JButton("Stop Send SMS");
        mStopButton.setBounds(186, 154, 125, 25);
        mStopButton.addActionListener(this);
        mMainFrame.getContentPane().add(mStopButton);
        /* Button: log export and log clear */
        mExportLogButton = new JButton("Export Log File");
        mExportLogButton.setBounds(26, 188, 125, 25);
        mMainFrame.getContentPane().add(mExportLogButton);
        mClearLogButton = new JButton("Clear Log File");
        mClearLogButton.setBounds(186, 188, 125, 25);
        mMainFrame.getContentPane().add(mClearLogButton);
        /* Change config button*/
        mLoadConfigButton = new JButton("Change Configuration");
        mLoadConfigButton.setBounds(26, 224, 285, 25);
        mLoadConfigButton.addActionListener(this);
        mMainFrame.getContentPane().add(mLoadConfigButton);
        /* Exit button */
        mExitButton = new JButton("Exit");
        mExitButton.setBounds(26, 260, 285, 25);
        mMainFrame.getContentPane().add(mExitButton);
    }
    /**
     * Init mMainFrame
     */
    private void initializeLoadConfig() {
        mLoadConfigFrame = new JFrame();
        mLoadConfigFrame.setBounds(100, 100, 360, 360);
        mLoadConfigFrame.setDefaultCloseOperation(JFrame.HIDE_ON_CLOSE);
        mLoadConfigFrame.getContentPane().setLayout(null);
        /*Ftp address */
        mFtpAddressTagLabel = new JLabel("Ftp Address: ");
        mFtpAddressTagLabel.setBounds(26, 32, 100, 16);
        mLoadConfigFrame.getContentPane().add(mFtpAddressTagLabel);
        mFtpAddressLabel = new JLabel("[...] ");
        mFtpAddressLabel.setBounds(136, 33, 175, 14);
        mLoadConfigFrame.getContentPane().add(mFtpAddressLabel);
        /* Fetch Config from SQL server checkbox*/
        mFetchConfigCheckBox = new JCheckBox("Load Configuration from SQL Server");
        mFetchConfigCheckBox.setBounds(26, 70, 250, 16);
        mFetchConfigCheckBox.addItemListener(this);
        mLoadConfigFrame.getContentPane().add(mFetchConfigCheckBox);
    }
    @Override
    public void actionPerformed(ActionEvent e){
        if (e.getSource() == this.mStartButton){
            Thread thread = new Thread(this.mJob);
            thread.start();
            this.mJobStatusLabel.setText("job started.");
        }
        if (e.getSource() == this.mStopButton){
            this.mJob.setActive(false);
            this.mJobStatusLabel.setText("job ended.");
        }
        if (e.getSource() == this.mLoadConfigButton){
            this.mLoadConfigFrame.setVisible(true);
        }
    }
    @Override
    public void itemStateChanged(ItemEvent e) {
        if (e.getSource() == this.mFetchConfigCheckBox){
            if (e.getStateChange() == 1){ //checked
                try {
                    SQLiteConfig config = new SQLiteConfig();
                    // config.setReadOnly(true);   
                    config.setSharedCache(true);
                    config.enableRecursiveTriggers(true);
                    SQLiteDataSource ds = new SQLiteDataSource(config); 
                    ds.setUrl("jdbc:sqlite:sms.db");
                    Connection con = ds.getConnection();
                    //ds.setServerName("sample.db");
                    String sql = "select * from sms_properties";
                    Statement stat = null;
                    ResultSet rs = null;
                    stat = con.createStatement();
                    rs = stat.executeQuery(sql);
                    while(rs.next()){
                        System.out.println(rs.getString("spftpadd")+"\t"+rs.getString("spftpprt"));
                    }
                }catch (Exception ex){
                    ex.printStackTrace();
                }
                this.mFtpAddressLabel.setText("to checked");
            }
            if (e.getStateChange() != 1){ //unchecked
                this.mFtpAddressLabel.setText("empty");
            }
        }
    }   
}

'@

$s158 = @'
JCheckbox x ItemListener
'@

$s187 = @'
import java.awt.event.ActionListener;
import java.awt.event.ActionEvent;
import java.awt.event.ItemListener;
import java.awt.event.ItemEvent;
import java.awt.EventQueue;
import java.awt.Font;
...
public class SendSMSUI implements ActionListener, ItemListener {
...
        /* Fetch Config from SQL server checkbox*/
        mFetchConfigCheckBox = new JCheckBox("Load Configuration from SQL Server");
        mFetchConfigCheckBox.setBounds(26, 70, 250, 16);
        mFetchConfigCheckBox.addItemListener(this);
        mLoadConfigFrame.getContentPane().add(mFetchConfigCheckBox);
...
    @Override
    public void itemStateChanged(ItemEvent e) {
        if (e.getSource() == this.mFetchConfigCheckBox){
            if (e.getStateChange() == 1){ //checked
                try {
                    SQLiteConfig config = new SQLiteConfig();
                    // config.setReadOnly(true);   
                    config.setSharedCache(true);
                    config.enableRecursiveTriggers(true);
                    SQLiteDataSource ds = new SQLiteDataSource(config); 
                    ds.setUrl("jdbc:sqlite:sms.db");
                    Connection con = ds.getConnection();
                    //ds.setServerName("sample.db");
                    String sql = "select * from sms_properties";
                    Statement stat = null;
                    ResultSet rs = null;
                    stat = con.createStatement();
                    rs = stat.executeQuery(sql);
                    while(rs.next()){
                        System.out.println(rs.getString("spftpadd")+"\t"+rs.getString("spftpprt"));
                    }
                }catch (Exception ex){
                    ex.printStackTrace();
                }
                this.mFtpAddressLabel.setText("to checked");
            }
            if (e.getStateChange() != 1){ //unchecked
                this.mFtpAddressLabel.setText("empty");
            }
        }
    }   
'@

$s159 = @'
Jbutton x ActionListener
'@

$s188 = @'
public class SendSMSUI implements ActionListener, ItemListener {
…
        /* Change config button*/
        mLoadConfigButton = new JButton("Change Configuration");
        mLoadConfigButton.setBounds(26, 224, 285, 25);
        mLoadConfigButton.addActionListener(this);
...
    @Override
    public void actionPerformed(ActionEvent e){
        if (e.getSource() == this.mStartButton){
            Thread thread = new Thread(this.mJob);
            thread.start();
            this.mJobStatusLabel.setText("job started.");
        }
        if (e.getSource() == this.mStopButton){
            this.mJob.setActive(false);
            this.mJobStatusLabel.setText("job ended.");
        }
        if (e.getSource() == this.mLoadConfigButton){
            this.mLoadConfigFrame.setVisible(true);
        }
    }
'@

$s160 = @'
JDBC
'@

$s161 = @'
SQLite POM
'@

$s177 = @'
    <!-- sqlite-jdbc -->
    <dependency>
      <groupId>org.xerial</groupId>
      <artifactId>sqlite-jdbc</artifactId>
      <version>3.7.2</version>
    </dependency>
'@

$s162 = @'
SQLite Connection
'@

$s194 = @'
                try {
                    SQLiteConfig config = new SQLiteConfig();
                    // config.setReadOnly(true);   
                    config.setSharedCache(true);
                    config.enableRecursiveTriggers(true);
                    SQLiteDataSource ds = new SQLiteDataSource(config); 
                    ds.setUrl("jdbc:sqlite:sms.db");
                    Connection con = ds.getConnection();
                    //ds.setServerName("sample.db");
                    String sql = "select * from sms_properties";
                    Statement stat = null;
                    ResultSet rs = null;
                    stat = con.createStatement();
                    rs = stat.executeQuery(sql);
                    while(rs.next()){
                        System.out.println(rs.getString("spftpadd")+"\t"+rs.getString("spftpprt"));
                    }
                }catch (Exception ex){
                    ex.printStackTrace();
                }
'@

$s163 = @'
POI
'@

$s184 = @'
POI Pom(Excel Part)
'@

$s183 = @'
    <!-- POI -->
    <dependency>
      <groupId>org.apache.poi</groupId>
      <artifactId>poi</artifactId>
      <version>3.16-beta1</version>
    </dependency>
    <!-- poi-ooxml -->
    <dependency>
        <groupId>org.apache.poi</groupId>
        <artifactId>poi-ooxml</artifactId>
        <version>3.16-beta1</version>
    </dependency>
'@

$s164 = @'
Basic
'@

$s185 = @'
package com.myles.bcm.poc;
import org.json.simple.JSONObject;
import org.json.simple.JSONArray;
import org.json.simple.parser.JSONParser;
import org.apache.poi.hssf.usermodel.HSSFWorkbook;
import org.apache.poi.ss.usermodel.*;
import org.apache.poi.hssf.util.CellReference;
import java.lang.Exception;
import java.io.FileOutputStream;
import java.io.FileInputStream;
import java.lang.Integer;
/**
 * The fetcher class
 *
 */
public class ExcelRateFetcher{
    final static int JSON_RESULT_FORMAT = 1;
    private String mResourceLocation; 
    private TargetObject mTargetObject;
    public ExcelRateFetcher(){
        super();
        this.mTargetObject = new TargetObject();
    }
    public void setResourceLocation(String location){
        this.mResourceLocation = location;
    }
    public void fetch(){
        try{
            Workbook workbook = WorkbookFactory.create(new FileInputStream(this.mResourceLocation));
            FormulaEvaluator evaluator = workbook.getCreationHelper().createFormulaEvaluator();
            Sheet sheet;
            CellReference cellReference;
            Row row;
            Cell cell;
            CellValue cellValue;
            /* Read Year */
            sheet = workbook.getSheetAt(0);
            cellReference = new CellReference("B2"); 
            row = sheet.getRow(cellReference.getRow());
            cell = row.getCell(cellReference.getCol()); 
            /*Debug*/ System.out.println("There is runtime exception with below line");
            cellValue = evaluator.evaluate(cell); 
            if (cellValue.getCellType() == Cell.CELL_TYPE_NUMERIC){
                this.mTargetObject.setYear(
                    new Double(cellValue.getNumberValue()).intValue()
                );
            }
            /* Read Month */
            sheet = workbook.getSheetAt(0);
            cellReference = new CellReference("B3"); 
            row = sheet.getRow(cellReference.getRow());
            cell = row.getCell(cellReference.getCol()); 
            cellValue = evaluator.evaluate(cell);
            if (cellValue.getCellType() == Cell.CELL_TYPE_NUMERIC){
                this.mTargetObject.setMonth(
                    new Double(cellValue.getNumberValue()).intValue()
                );
            }
            /* Read Day */
            sheet = workbook.getSheetAt(0);
            cellReference = new CellReference("B4"); 
            row = sheet.getRow(cellReference.getRow());
            cell = row.getCell(cellReference.getCol()); 
            cellValue = evaluator.evaluate(cell);
            if (cellValue.getCellType() == Cell.CELL_TYPE_NUMERIC){
                this.mTargetObject.setDay(
                    new Double(cellValue.getNumberValue()).intValue()
                );
            }
            /* Read Hour */
            sheet = workbook.getSheetAt(0);
            cellReference = new CellReference("B5"); 
            row = sheet.getRow(cellReference.getRow());
            cell = row.getCell(cellReference.getCol()); 
            cellValue = evaluator.evaluate(cell);
            if (cellValue.getCellType() == Cell.CELL_TYPE_NUMERIC){
                this.mTargetObject.setHour(
                    new Double(cellValue.getNumberValue()).intValue()
                );
            }
            /* Read Minute */
            sheet = workbook.getSheetAt(0);
            cellReference = new CellReference("B6"); 
            row = sheet.getRow(cellReference.getRow());
            cell = row.getCell(cellReference.getCol()); 
            cellValue = evaluator.evaluate(cell);
            if (cellValue.getCellType() == Cell.CELL_TYPE_NUMERIC){
                this.mTargetObject.setMinute(
                    new Double(cellValue.getNumberValue()).intValue()
                );
            }
            /* Read Second */
            sheet = workbook.getSheetAt(0);
            cellReference = new CellReference("B7"); 
            row = sheet.getRow(cellReference.getRow());
            cell = row.getCell(cellReference.getCol()); 
            cellValue = evaluator.evaluate(cell);
            if (cellValue.getCellType() == Cell.CELL_TYPE_NUMERIC){
                this.mTargetObject.setSecond(
                    new Double(cellValue.getNumberValue()).intValue()
                );
            }
        }catch(Exception e){
            System.out.println(e);
        }    
    }
    public Object getRate(int format){
        if (format == ExcelRateFetcher.JSON_RESULT_FORMAT){
            try{
                JSONObject jsonObject = new JSONObject();
                jsonObject.put(
                    TargetObject.YEAR_FIELDNAME, 
                    new Integer(this.mTargetObject.getYear()).toString()
                );
                jsonObject.put(
                    TargetObject.MONTH_FIELDNAME, 
                    new Integer(this.mTargetObject.getMonth()).toString()
                );
                jsonObject.put(
                    TargetObject.DAY_FIELDNAME, 
                    new Integer(this.mTargetObject.getDay()).toString()
                );
                jsonObject.put(
                    TargetObject.HOUR_FIELDNAME,
                    new Integer(this.mTargetObject.getHour()).toString()
                );
                jsonObject.put(
                    TargetObject.MINUTE_FIELDNAME, 
                    new Integer(this.mTargetObject.getMinute()).toString()
                );
                jsonObject.put(
                    TargetObject.SECOND_FIELDNAME, 
                    new Integer(this.mTargetObject.getSecond()).toString()
                );
                return jsonObject;
            }catch(Exception pe){
                System.out.println(pe);
                return null;
            }
        }else{
            System.out.println("Output format not supported.");
            return null;
        }
    }
}

'@

$s165 = @'
Log4j
'@

$s166 = @'
Log4j POM
'@

$s182 = @'
    <!-- log4j -->
    <dependency>
      <groupId>log4j</groupId>
      <artifactId>log4j</artifactId>
      <version>1.2.17</version>
    </dependency>
'@

$s167 = @'
Basic 
'@

$s190 = @'
public class SendSMSJob extends Thread {
    final static Logger mLogger = Logger.getLogger(SendSMSJob.class);
    private boolean mIsActive;
    private List<FileManipulator> mProcessChain;
    private String FTP_ADDRESS;
    private int FTP_PORT;
    private String FTP_USER;
    private String FTP_PASSWORD;
    private String FTP_FOLDER;
    private String SMS_FILE;
    private String SMS_FOLDER;
    private String BKUP_FOLDER;
    private String CUR_DIR;
    private String INMSG;
    private String OUTMSG;
    private int LOOP_INTERVAL;
    private String LOG_PROP;   
    public SendSMSJob(){
        super();
        this.mProcessChain = new ArrayList<FileManipulator>();
        this.mIsActive = false;
    }
    public SendSMSJob(String propertiesPath) {
        super();
        this.mProcessChain = new ArrayList<FileManipulator>();
        this.mIsActive = false;
        this.setProperties(propertiesPath);        
    }
    public void setProperties(String propertiesPath){
        Properties prop = new Properties();
        InputStream input = null;
        try {
            input = new FileInputStream(propertiesPath);
            prop.load(input);
            FTP_ADDRESS = prop.getProperty("FTP_ADDRESS");
            FTP_PORT = Integer.parseInt(prop.getProperty("FTP_PORT"));
            FTP_USER = prop.getProperty("FTP_USER");
            FTP_PASSWORD = prop.getProperty("FTP_PASSWORD");
            FTP_FOLDER = prop.getProperty("FTP_FOLDER");
            SMS_FILE = prop.getProperty("SMS_FILE");
            SMS_FOLDER = prop.getProperty("SMS_FOLDER");
            BKUP_FOLDER = prop.getProperty("BKUP_FOLDER");
            CUR_DIR = System.getProperty("user.dir");
            INMSG = prop.getProperty("INMSG");
            OUTMSG = prop.getProperty("OUTMSG");
            LOOP_INTERVAL = Integer.parseInt(prop.getProperty("LOOP_INTERVAL"));
            LOG_PROP = prop.getProperty("LOG_PROP");
        } catch (IOException ex) {
            ex.printStackTrace();
        }
        /* Set up MessageRegisterLogProxy */
        MessageRegisterLogProxy messageRegisterLogProxy = new MessageRegisterLogProxy();
        /* Set up MessageFtpUploaderLogProxy */
        MessageFtpUploaderLogProxy messageFtpUploaderLogProxy = new MessageFtpUploaderLogProxy();
        messageFtpUploaderLogProxy.setFtpAddress(FTP_ADDRESS);
        messageFtpUploaderLogProxy.setFtpPort(FTP_PORT);
        messageFtpUploaderLogProxy.setFtpUser(FTP_USER);
        messageFtpUploaderLogProxy.setFtpPassword(FTP_PASSWORD);
        messageFtpUploaderLogProxy.setFtpFolder(FTP_FOLDER);
        /* Set up MessageBackuperLogProxy */
        MessageBackuperLogProxy messageBackuperLogProxy = new MessageBackuperLogProxy();
        messageBackuperLogProxy.setPath(BKUP_FOLDER);
        this.mProcessChain.add(messageRegisterLogProxy);
        this.mProcessChain.add(messageFtpUploaderLogProxy);
        this.mProcessChain.add(messageBackuperLogProxy);
        /* Ddynamic configuration of log setting */
        PropertyConfigurator.configure(LOG_PROP);
    }
    public boolean isActive(){
        return this.mIsActive;
    }
    public void setActive(boolean b){
        this.mIsActive = b;
    }
    @Override
    public void run(){
        mLogger.debug("Job Starts.");
        this.mIsActive = true;
        while(this.isActive()) {
            try {
                //File operations
                File targetFoler = new File(SMS_FOLDER);
                if (targetFoler.exists() && targetFoler.isDirectory()){
                    for (File f : targetFoler.listFiles()){
                        // System.out.println("Found file(s):" + f);
                        // mLogger.info("Found file(s):" + f);
                        for (FileManipulator fm : mProcessChain){
                            fm.setFile(f);
                            fm.manipulate();
                            if (!fm.isSuccess()){
                                break;
                            }
                        }
                    }
                }
                // 暫停目前的執行緒5秒
                Thread.sleep(LOOP_INTERVAL);
                mLogger.debug("Job loops.");
            } catch(InterruptedException e) {
                e.printStackTrace();
            }
        }
        mLogger.debug("Job Ends.");
    }
    public static void main(String[] args) {
        SendSMSJob job = new SendSMSJob(args[0]);
        job.run(); //not using thread, start using the current thread, terminated by Ctrl C in bash
    }
}
'@

$s189 = @'
Override Properties at runtime
'@

$s193 = @'
PropertyConfigurator.configure(LOG_PROP);   //LOG_PROP is a string of path
'@

$s168 = @'
Properties
'@

$s192 = @'
Basic (load properties_)
'@

$s191 = @'
    public void setProperties(String propertiesPath){
        Properties prop = new Properties();
        InputStream input = null;
        try {
            input = new FileInputStream(propertiesPath);
            prop.load(input);
            FTP_ADDRESS = prop.getProperty("FTP_ADDRESS");
            FTP_PORT = Integer.parseInt(prop.getProperty("FTP_PORT"));
            FTP_USER = prop.getProperty("FTP_USER");
            FTP_PASSWORD = prop.getProperty("FTP_PASSWORD");
            FTP_FOLDER = prop.getProperty("FTP_FOLDER");
            SMS_FILE = prop.getProperty("SMS_FILE");
            SMS_FOLDER = prop.getProperty("SMS_FOLDER");
            BKUP_FOLDER = prop.getProperty("BKUP_FOLDER");
            CUR_DIR = System.getProperty("user.dir");
            INMSG = prop.getProperty("INMSG");
            OUTMSG = prop.getProperty("OUTMSG");
            LOOP_INTERVAL = Integer.parseInt(prop.getProperty("LOOP_INTERVAL"));
            LOG_PROP = prop.getProperty("LOG_PROP");
        } catch (IOException ex) {
            ex.printStackTrace();
        }
        
'@

$s169 = @'
Junit
'@

$s170 = @'
Junit POM (Annotation @Before is only support after 4)
'@

$s173 = @'
    <!-- JUnit 4 -->
    <dependency>
      <groupId>junit</groupId>
      <artifactId>junit</artifactId>
      <version>4.11</version>
    </dependency>
'@

$s171 = @'
JUnit
'@

$s172 = @'
Basic usage
'@

$s195 = @'
package com.bcm.app.core;
import static org.junit.Assert.assertEquals;
import org.junit.*;
import java.io.File;
public class MessageRegisterTest{
    private FileManipulator messageRegister;
    @Before
    public void createInstance(){
        // Reset the object e
        messageRegister = new MessageRegister();
    }
    @Test
    public void testGetFileEmpty() {
        assertEquals(messageRegister.getFile(), null);
    }
    @Test 
    public void testSetNullFile(){
        try{
            messageRegister.setFile(null);
            assertEquals(true, true);  //Exception test: shd not throws exception
        }catch(Exception e ){
            assertEquals(true, false);  //Exception test
        }
    }
    @Test
    public void testGetFile() {
        File f = new File("test");
        messageRegister.setFile(f);
        assertEquals(messageRegister.getFile(), f);
    }
    @Test 
    public void testFailBeforeManipulate(){
        assertEquals(messageRegister.isSuccess(), false);
    }
    @Test 
    public void testFailAfterManipulateWithoutSetFile(){
        messageRegister.manipulate();
        assertEquals(messageRegister.isSuccess(), false);
    }
    @Test 
    public void testFailAfterManipulateWithNonExistingFile(){
        File f = new File("test");
        messageRegister.setFile(f);
        messageRegister.manipulate();
        assertEquals(messageRegister.isSuccess(), false);
    }
    @Test 
    public void testSuccessAfterManipulateWithExistingFile(){
        try{
            String path = System.getProperty("user.dir");
            File f = new File(path + "\\test.txt");
            if(!f.exists()){
                f.createNewFile();
            }
            messageRegister.setFile(f);        
            messageRegister.manipulate();
            assertEquals(messageRegister.isSuccess(), true);
            f.delete();
        }catch (Exception e){
            e.printStackTrace();
        }
    }
    @Test 
    public void testIsSuccessResetWhenSetExistingFileThenSetNonExistingFile(){
        try{
            String path = System.getProperty("user.dir");
            File f1 = new File(path + "\\test1.txt");
            if(!f1.exists()){
                f1.createNewFile();
            }
            File f2 = new File(path + "\\test2.txt");
            if(f2.exists()){
                f2.delete();
            }
            messageRegister.setFile(f1);        
            messageRegister.manipulate();
            boolean firstTimeStatus = messageRegister.isSuccess();
            messageRegister.setFile(f2);        
            boolean secondTimeStatus = messageRegister.isSuccess();
            assertEquals(firstTimeStatus && !secondTimeStatus, true);
            f1.delete();
        }catch (Exception e){
            e.printStackTrace();
        }
    }
}

'@

$s174 = @'
Apache Common Net
'@

$s175 = @'
FTP Basic Udage
'@

$s197 = @'
package com.bcm.app.core;
import java.io.*;
import org.apache.commons.net.ftp.FTPClient;
import org.apache.commons.net.ftp.FTPReply;
public class MessageFtpUploader implements FileManipulator{
    private File mFile;
    private boolean mIsSuccess = false;
    private String mFtpAddress;
    private int mFtpPort;
    private String mFtpUser;
    private String mFtpPassword;
    private String mFtpFolder;
    /* --- Setters ---*/
    public void setFtpAddress(String address){
        this.mFtpAddress = address;
    }
    public void setFtpPort(int port){
        this.mFtpPort = port;
    }
    public void setFtpUser(String user){
        this.mFtpUser = user;
    }
    public void setFtpPassword(String password){
        this.mFtpPassword = password;
    }
    public void setFtpFolder(String folder){
        this.mFtpFolder = folder;
    }
    public String getFtpAddress(){
        return this.mFtpAddress;
    }
    public int getFtpPort(){
        return this.mFtpPort;
    }
    public String getFtpUser(){
        return this.mFtpUser;
    }
    public String getFtpPassword(){
        return this.mFtpPassword;
    }
    public String getFtpFolder(){
        return this.mFtpFolder;
    }
    /* --- interface override ---*/
    @Override
    public void setFile(File file){
        this.mFile = file;
        this.mIsSuccess = false;
    };
    @Override
    public File getFile(){
        return this.mFile;
    };    
    /**
     * Method manipulate in MessageUploader will check the 
     * existence of the file set to the object
     */
    @Override
    public void manipulate(){
        try{
            this.mIsSuccess = false;
            FTPClient ftpClient = new FTPClient();
            ftpClient.connect(this.getFtpAddress(), this.getFtpPort());
            int replyCode = ftpClient.getReplyCode();
            if (!FTPReply.isPositiveCompletion(replyCode)) {
                System.out.println("Operation failed. Server reply code: " + replyCode);
                this.mIsSuccess = false;
                return ;
            }
            boolean success = ftpClient.login(this.getFtpUser(), this.getFtpPassword());
            if (!success){            
                System.out.println("Wrong ftp settings, test skiped.");
                this.mIsSuccess = false;
                return ;
            }
            ftpClient.changeWorkingDirectory(this.getFtpFolder());
            if (!FTPReply.isPositiveCompletion(replyCode)){
                System.out.println("Cannot find target ftp folder, test skipped. ");
                this.mIsSuccess = false;
                return ; 
            }
            InputStream input = new FileInputStream(this.getFile());  //fileName includes filetype
            ftpClient.appendFile(this.getFile().getName(), input);
            input.close();
            ftpClient.logout();
            ftpClient.disconnect();
            this.mIsSuccess = true;
        }catch (Exception e){
            e.printStackTrace();
        }
    }
    @Override
    public boolean isSuccess(){
        return mIsSuccess;
    }   
}
'@

$s176 = @'
POM
'@

$s178 = @'
    <!-- Apache common net -->
    <dependency>
      <groupId>commons-net</groupId>
      <artifactId>commons-net</artifactId>
      <version>3.3</version>
    </dependency>
'@

$s179 = @'
Apache common io
'@

$s181 = @'
    <!--Apache common io-->
    <dependency>
      <groupId>commons-io</groupId>
      <artifactId>commons-io</artifactId>
      <version>2.5</version>
    </dependency>
'@

$s180 = @'
Bsaic (Super Convenient, like bash util)
'@

$s196 = @'
import org.apache.commons.io.FileUtils;
…
            File copiedFile = new File("copiedFile.txt");
            FileUtils.touch(copiedFile);
            FileUtils.copyFile(f, copiedFile);
…
            assertEquals(FileUtils.contentEquals(copiedFile, backupFile), true);  //for backup isPerformed
'@


$ws1.Rows.Item(40).RowHeight = 33
$ws1.Cells.Item(40,1).Value2 = $s156
$ws1.Cells.Item(40,2).Value2 = $s157
$ws1.Cells.Item(40,3).Value2 = $s186

$ws1.Rows.Item(41).RowHeight = 33
$ws1.Cells.Item(41,1).Value2 = $s156
$ws1.Cells.Item(41,2).Value2 = $s158
$ws1.Cells.Item(41,3).Value2 = $s187

$ws1.Rows.Item(42).RowHeight = 33
$ws1.Cells.Item(42,1).Value2 = $s156
$ws1.Cells.Item(42,2).Value2 = $s159
$ws1.Cells.Item(42,3).Value2 = $s188

$ws1.Rows.Item(43).RowHeight = 33
$ws1.Cells.Item(43,1).Value2 = $s160
$ws1.Cells.Item(43,2).Value2 = $s161
$ws1.Cells.Item(43,3).Value2 = $s177

$ws1.Rows.Item(44).RowHeight = 33
$ws1.Cells.Item(44,1).Value2 = $s160
$ws1.Cells.Item(44,2).Value2 = $s162
$ws1.Cells.Item(44,3).Value2 = $s194

$ws1.Rows.Item(45).RowHeight = 33
$ws1.Cells.Item(45,1).Value2 = $s163
$ws1.Cells.Item(45,2).Value2 = $s184
$ws1.Cells.Item(45,3).Value2 = $s183

$ws1.Rows.Item(46).RowHeight = 33
$ws1.Cells.Item(46,1).Value2 = $s163
$ws1.Cells.Item(46,2).Value2 = $s164
$ws1.Cells.Item(46,3).Value2 = $s185

$ws1.Rows.Item(47).RowHeight = 33
$ws1.Cells.Item(47,1).Value2 = $s165
$ws1.Cells.Item(47,2).Value2 = $s166
$ws1.Cells.Item(47,3).Value2 = $s182

$ws1.Rows.Item(48).RowHeight = 33
$ws1.Cells.Item(48,1).Value2 = $s165
$ws1.Cells.Item(48,2).Value2 = $s167
$ws1.Cells.Item(48,3).Value2 = $s190

$ws1.Rows.Item(49).RowHeight = 33
$ws1.Cells.Item(49,1).Value2 = $s165
$ws1.Cells.Item(49,2).Value2 = $s189
$ws1.Cells.Item(49,3).Value2 = $s193

$ws1.Rows.Item(50).RowHeight = 33
$ws1.Cells.Item(50,1).Value2 = $s168
$ws1.Cells.Item(50,2).Value2 = $s192
$ws1.Cells.Item(50,3).Value2 = $s191

$ws1.Rows.Item(51).RowHeight = 33
$ws1.Cells.Item(51,1).Value2 = $s169
$ws1.Cells.Item(51,2).Value2 = $s170
$ws1.Cells.Item(51,3).Value2 = $s173

$ws1.Rows.Item(52).RowHeight = 33
$ws1.Cells.Item(52,1).Value2 = $s171
$ws1.Cells.Item(52,2).Value2 = $s172
$ws1.Cells.Item(52,3).Value2 = $s195

$ws1.Rows.Item(53).RowHeight = 33
$ws1.Cells.Item(53,1).Value2 = $s174
$ws1.Cells.Item(53,2).Value2 = $s175
$ws1.Cells.Item(53,3).Value2 = $s197

$ws1.Rows.Item(54).RowHeight = 33
$ws1.Cells.Item(54,1).Value2 = $s174
$ws1.Cells.Item(54,2).Value2 = $s176
$ws1.Cells.Item(54,3).Value2 = $s178

$ws1.Rows.Item(55).RowHeight = 33
$ws1.Cells.Item(55,1).Value2 = $s179
$ws1.Cells.Item(55,2).Value2 = $s176
$ws1.Cells.Item(55,3).Value2 = $s181

$ws1.Rows.Item(56).RowHeight = 33
$ws1.Cells.Item(56,1).Value2 = $s179
$ws1.Cells.Item(56,2).Value2 = $s180
$ws1.Cells.Item(56,3).Value2 = $s196

# sheet1 view state
$ws1.Application.ActiveWindow.ScrollRow = 45
$ws1.Range("C51").Select()

# ---------- Sheet3 ("the source"): insert 5 rows above row 3, add new entries ----------
$ws3 = $wb.Worksheets.Item("the source")
$ws3.Rows.Item(3).Resize(5).Insert()
$ws3.Columns.Item(3).ColumnWidth = 39.28515625

$t150 = @'
C:\Users\BI77\Documents\identical_ref\jdk1.8_src
'@

$t151 = @'
java.util.Collection
'@

$t152 = @'
java.lang.Boolean
'@

$t153 = @'
synchronized method in Thread
'@

$t154 = @'
java.util.Properties
java.util.Hashtable
java.util.Enumeration<T>
'@

$t155 = @'
一個簡單的接口, hasMoreElement和nextElement, 在Properties看到有所應用, 當要遍歷HashTable時, 將h,key()賦給一個臨時的Enumeration 類(這key大有文章, 是由HashTable中的inner class Enumerator 實現的), 再以for 的三段式利用hasMoreElement()和nextElement()來Iterate
'@


$ws3.Cells.Item(2,1).Value2 = $t150

$ws3.Rows.Item(6).RowHeight = 60
$ws3.Cells.Item(6,1).Value2 = 42816
$ws3.Cells.Item(6,1).NumberFormat = "m/d/yyyy"
$ws3.Cells.Item(6,2).Value2 = $t155
$ws3.Cells.Item(6,3).Value2 = $t154

$ws3.Cells.Item(7,1).Value2 = 42814
$ws3.Cells.Item(7,1).NumberFormat = "m/d/yyyy"
$ws3.Cells.Item(7,2).Value2 = $t153

$ws3.Cells.Item(8,3).Value2 = $t151
$ws3.Cells.Item(9,3).Value2 = $t152

$ws3.Range("A33").Select()

